# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (Home) totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 531
$wsOff.Range("C2").Value = 368
$wsOff.Range("D2").Value = 113
$wsOff.Range("E2").Value = 57
$wsOff.Range("F2").Value = 9

# DEF sheet - row 2 (Home) totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 457
$wsDef.Range("C2").Value = 314
$wsDef.Range("D2").Value = 106
$wsDef.Range("E2").Value = 45
$wsDef.Range("F2").Value = 7
